$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace "proton" with "p" in the target column (G2:G13) for every data row.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 7).Value = "p"
}

# Make the header row (A1:K1) bold and centered.
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
